# Add new columns I (I0) and J (IF) to the worksheet, with header styling
# matching the existing header row, and fill in the data values for rows 2-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 onto the new header cells I1:J1, then set
# their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I0 / IF columns, keyed by row number.
$i0Values = @{
    2 = 1
    3 = 1
    4 = 1
    5 = 1
    6 = 1
    7 = 1
    8 = 1
    9 = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 4
    14 = 6
    15 = 5
}

$ifValues = @{
    2 = 6
    3 = 6
    4 = 6
    5 = 4
    6 = 7
    7 = 4
    8 = 5
    9 = 4
    10 = 6
    11 = 6
    12 = 5
    13 = 6
    14 = 7
    15 = 6
}

foreach ($r in 2..15) {
    $ws.Cells.Item($r, 9).Value = $i0Values[$r]
    $ws.Cells.Item($r, 10).Value = $ifValues[$r]
}
